# Adding the changes we made on may 9th
# Insert 5 new rows of accelerometer data at the top (rows 2-6, pushing
# existing data down) and append 5 new rows of data at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 5 rows before the current row 2 (existing data shifts down) ---
$insertRange = $ws.Range("A2:C6")
$insertRange.Insert()
# Excel copies the formatting of the row above (the bold header row) onto
# newly inserted rows; clear it so the new data rows stay unstyled, matching
# the rest of the plain numeric data rows.
$insertRange.ClearFormats()

$topData = @(
    @(-2.278594136238097, 10.26830673217773, -1.992950439453123),
    @(-1.48284924030304, 9.543427348136902, -0.6066013872623449),
    @(-2.347809791564941, 9.136160850524902, -1.348505258560181),
    @(-2.276926577091217, 9.172239780426027, -1.260899052023888),
    @(-3.074454665184022, 9.003937959671021, -1.032875627279282)
)

for ($i = 0; $i -lt $topData.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $topData[$i][0]
    $ws.Cells.Item($row, 2).Value = $topData[$i][1]
    $ws.Cells.Item($row, 3).Value = $topData[$i][2]
}

# --- Append 5 rows after the last existing row (now row 26) ---
$bottomData = @(
    @(3.198545038700105, 8.029186010360718, 1.33356249332428),
    @(3.5916051864624, 8.197292327880859, 1.436704874038696),
    @(3.165686726570128, 8.038311719894409, 1.441773623228073),
    @(3.061142683029175, 7.987302541732788, 1.411497831344604),
    @(3.059285700321198, 8.002557039260864, 1.58986583352089)
)

for ($i = 0; $i -lt $bottomData.Length; $i++) {
    $row = 27 + $i
    $ws.Cells.Item($row, 1).Value = $bottomData[$i][0]
    $ws.Cells.Item($row, 2).Value = $bottomData[$i][1]
    $ws.Cells.Item($row, 3).Value = $bottomData[$i][2]
}
